$wb = $excel.ActiveWorkbook

function Set-OrClear($ws, $cellRef, $val) {
    if ($null -eq $val) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $val
    }
}

$ws = $wb.Worksheets.Item("ALC")
$changes = @{
    "H2" = 1424.6316
    "I2" = 350
    "J2" = 2051.5
    "K2" = 350
    "L2" = 2051.5
    "M2" = -237
    "N2" = -2277.5
    "H4" = 394.4
    "J4" = 681.6667
    "L4" = 681.6667
    "N4" = -909.6667
    "H6" = 273.6154
    "I6" = 279.72726
    "K6" = 839.18178
    "M6" = -727.18178
    "H9" = 580.0909
    "I9" = 395.57144
    "J9" = 903
    "K9" = 395.57144
    "L9" = 903
    "M9" = -226.57144
    "N9" = -1241
    "H15" = 2219.25
    "I15" = 2219.25
    "K15" = 6657.75
    "M15" = -6488.75
    "H17" = 551.9268
    "J17" = 551.9268
    "L17" = 1655.7804
    "N17" = -1991.7804
    "H18" = 1798.4286
    "I18" = 1931.6666
    "J18" = 999
    "K18" = 1931.6666
    "L18" = 999
    "M18" = -1647.6666
    "N18" = -1567
    "H33" = 7012.067
    "I33" = 8489.333000000001
    "K33" = 8489.333000000001
    "M33" = -8260.333000000001
    "H64" = 5426.6665
    "J64" = 5500
    "L64" = 5500
    "N64" = -5996
    "H67" = 5426.6665
    "J67" = 5500
    "L67" = 5500
    "N67" = -7216
    "H69" = 8172
    "I69" = 7035.75
    "J69" = 9990
    "K69" = 21107.25
    "L69" = 29970
    "M69" = -20233.25
    "N69" = -31718
    "H72" = 8172
    "I72" = 7035.75
    "J72" = 9990
    "K72" = 63321.75
    "L72" = 89910
    "M72" = -58953.75
    "N72" = -98646
    "H98" = 1350.1428
    "I98" = 1321.1724
    "K98" = 1321.1724
    "M98" = 176.8276000000001
    "H122" = 1350.1428
    "I122" = 1321.1724
    "K122" = 3963.5172
    "M122" = -1513.5172
    "H128" = 75000
    "J128" = 0
    "L128" = 0
    "N128" = $null
    "H132" = 6525.0835
    "I132" = 6836.4546
    "K132" = 20509.3638
    "M132" = -17979.3638
    "H135" = 494.66666
    "I135" = 494.66666
    "K135" = 4451.99994
    "M135" = -1916.99994
    "H137" = 1732.7407
    "I137" = 1363.3636
    "K137" = 4090.0908
    "M137" = -1540.0908
    "H138" = 5882.8
    "I138" = 6926.1113
    "K138" = 20778.3339
    "M138" = -15638.3339
}
foreach ($cellRef in $changes.Keys) {
    Set-OrClear $ws $cellRef $changes[$cellRef]
}

$ws = $wb.Worksheets.Item("ARM")
$changes = @{
    "H45" = 2944.524
    "I45" = 1549.75
    "J45" = 4804.222
    "K45" = 1549.75
    "L45" = 4804.222
    "M45" = -1172.75
    "N45" = -5558.222
    "H70" = 74999
    "J70" = 74999
    "L70" = 74999
    "N70" = -75539
    "H73" = 74999
    "J73" = 74999
    "L73" = 74999
    "N73" = -76871
    "H128" = 0
    "J128" = 0
    "L128" = 0
    "N128" = $null
}
foreach ($cellRef in $changes.Keys) {
    Set-OrClear $ws $cellRef $changes[$cellRef]
}

$ws = $wb.Worksheets.Item("BSM")
$changes = @{
    "H86" = 2999.258
    "I86" = 1489.8823
    "K86" = 1489.8823
    "M86" = -366.8823
    "H89" = 2999.258
    "I89" = 1489.8823
    "K89" = 7449.4115
    "M89" = -1833.4115
    "H134" = 2770.3242
    "I134" = 1979.4839
    "K134" = 5938.4517
    "M134" = -3403.4517
    "H139" = 69997.39999999999
    "J139" = 69996.75
    "L139" = 69996.75
    "N139" = -80276.75
    "H141" = 79997
    "J141" = 79997
    "L141" = 79997
    "N141" = -90357
}
foreach ($cellRef in $changes.Keys) {
    Set-OrClear $ws $cellRef $changes[$cellRef]
}

$ws = $wb.Worksheets.Item("CRP")
$changes = @{
    "H7" = 244.625
    "I7" = 89.111115
    "K7" = 89.111115
    "M7" = 23.888885
    "H41" = 1249.75
    "I41" = 1249.75
    "K41" = 1249.75
    "M41" = -821.75
    "H100" = 67499.25
    "J100" = 76666
    "L100" = 76666
    "N100" = -78830
    "H122" = 2199.6667
    "I122" = 2199.6667
    "K122" = 6599.000100000001
    "M122" = -4149.000100000001
    "H137" = 75354
    "I137" = 50709
    "K137" = 50709
    "M137" = -45609
}
foreach ($cellRef in $changes.Keys) {
    Set-OrClear $ws $cellRef $changes[$cellRef]
}

$ws = $wb.Worksheets.Item("CUL")
$changes = @{
    "H35" = 706
    "I35" = 59
    "J35" = 2000
    "K35" = 177
    "L35" = 6000
    "M35" = 111
    "N35" = -6576
    "H39" = 1465.6666
    "I39" = 699
    "K39" = 2097
    "M39" = -1803
    "H113" = 570
    "I113" = 773
    "J113" = 512
    "K113" = 2319
    "L113" = 1536
    "M113" = -149
    "N113" = -5876
    "H140" = 2616.182
    "I140" = 2345.1875
    "K140" = 7035.5625
    "M140" = -1855.5625
}
foreach ($cellRef in $changes.Keys) {
    Set-OrClear $ws $cellRef $changes[$cellRef]
}

$ws = $wb.Worksheets.Item("GSM")
$changes = @{
    "H126" = 5583.5713
    "I126" = 4909.875
    "K126" = 14729.625
    "M126" = -12259.625
    "H132" = 50136.332
    "I132" = 55065.74
    "J132" = 3307
    "K132" = 165197.22
    "L132" = 9921
    "M132" = -162667.22
    "N132" = -14981
}
foreach ($cellRef in $changes.Keys) {
    Set-OrClear $ws $cellRef $changes[$cellRef]
}

$ws = $wb.Worksheets.Item("LTW")
$changes = @{
    "H55" = 235.77777
    "I55" = 262.57144
    "J55" = 142
    "K55" = 262.57144
    "L55" = 142
    "M55" = -89.57144
    "N55" = -488
    "H82" = 2325.5557
    "I82" = 1962.3636
    "J82" = 2896.2856
    "K82" = 1962.3636
    "L82" = 2896.2856
    "M82" = -1601.3636
    "N82" = -3618.2856
    "H85" = 2325.5557
    "I85" = 1962.3636
    "J85" = 2896.2856
    "K85" = 1962.3636
    "L85" = 2896.2856
    "M85" = -714.3635999999999
    "N85" = -5392.2856
    "H100" = 4534.5
    "I100" = 4319.3335
    "K100" = 4319.3335
    "M100" = -3778.3335
    "H132" = 40472.184
    "I132" = 48830.152
    "J132" = 9428.286
    "K132" = 146490.456
    "L132" = 28284.858
    "M132" = -143960.456
    "N132" = -33344.858
}
foreach ($cellRef in $changes.Keys) {
    Set-OrClear $ws $cellRef $changes[$cellRef]
}

$ws = $wb.Worksheets.Item("WVR")
$changes = @{
    "H43" = 49999.5
    "J43" = 49999.5
    "L43" = 49999.5
    "N43" = -50297.5
    "H56" = 0
    "J56" = 0
    "L56" = 0
    "N56" = $null
    "H136" = 3335.6
    "I136" = 2966.7144
    "K136" = 8900.143199999999
    "M136" = -6350.143199999999
}
foreach ($cellRef in $changes.Keys) {
    Set-OrClear $ws $cellRef $changes[$cellRef]
}
